# Update view-count figures for two 漫展 (anime convention) events across
# the "展览" sheet and the "全部类型" sheet, matching the regenerated
# gh-pages data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1067
$wsExhibit.Range("F5").Value = 563

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1067
$wsAll.Range("F6").Value = 563
